# JUL update, lab foundation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in the existing "British Fellowship awarding body" entry (row 17, column B)
$ws.Range("B17").Value = "British Fellowship awarding body (name witheld)"

# Add two new rows describing additional reviewing activity
$ws.Range("A18").Value = "Article review"
$ws.Range("B18").Value = "European Journal of Immunology"
$ws.Range("C18").Value = 2024

$ws.Range("A19").Value = "Poster session reviewer"
$ws.Range("B19").Value = "European Congress of Immunology"
$ws.Range("C19").Value = 2024

# Update the selected cell to match the author's final cursor position
$ws.Range("B13").Select()
